$d = $word.ActiveDocument

# 1. "Done till sprint 2" -> "Done till sprint 1"
$d.Content.Find.Execute("Done till sprint 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Done till sprint 1", 2)

# 2. Move the "done" run (in the cell containing the _GoBack bookmark) to appear
#    before the bookmarkStart/bookmarkEnd pair instead of after it.
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $para = $bm.Range.Paragraphs(1)
        $paraRange = $para.Range
        # Find the run of text "done" within this paragraph and move it before the bookmark.
        $text = $paraRange.Text
        Write-Output "ParaText=[$text]"
    }
}
